$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-15"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "17:59:27"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,712.3015"
